$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.72
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 4.65
$ws.Range("J3").Value = 2.25
$ws.Range("L3").Value = 4.8
$ws.Range("Q3").Value = 1.93
$ws.Range("W3").Value = 6.3
$ws.Range("X3").Value = 7.6
$ws.Range("Z3").Value = 13.5
$ws.Range("AH3").Value = 11.75
$ws.Range("AI3").Value = 26
$ws.Range("AJ3").Value = 15
$ws.Range("AK3").Value = 80
$ws.Range("AP3").Value = 17.5
$ws.Range("AQ3").Value = 28
$ws.Range("AS3").Value = 250
$ws.Range("AW3").Value = 6.3

# Row 4
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 2.6
$ws.Range("K4").Value = 2.3
$ws.Range("L4").Value = 3.75
$ws.Range("U4").Value = 1.62
$ws.Range("V4").Value = 2.2
$ws.Range("X4").Value = 11
$ws.Range("Z4").Value = 19
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 21
$ws.Range("AH4").Value = 13
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 12
$ws.Range("AK4").Value = 41
$ws.Range("AL4").Value = 26
$ws.Range("AO4").Value = 11
$ws.Range("AQ4").Value = 34
$ws.Range("AS4").Value = 101
$ws.Range("AW4").Value = 5.5
$ws.Range("AX4").Value = 19

# Row 7
$ws.Range("G7").Value = 1.67
$ws.Range("I7").Value = 5.75
$ws.Range("J7").Value = 2.4
$ws.Range("L7").Value = 7
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6
$ws.Range("Q7").Value = 2.88
$ws.Range("R7").Value = 1.4
$ws.Range("X7").Value = 6
$ws.Range("Z7").Value = 12
$ws.Range("AD7").Value = 7.5
$ws.Range("AE7").Value = 29
$ws.Range("AW7").Value = 7.5
$ws.Range("AZ7").Value = 201
$ws.Range("BA7").Value = 301

# Row 9
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 6.5

# Row 10
$ws.Range("G10").Value = 2.75
$ws.Range("I10").Value = 2.55
$ws.Range("J10").Value = 3.5
$ws.Range("L10").Value = 3.25
$ws.Range("N10").Value = 9
$ws.Range("W10").Value = 8
$ws.Range("AA10").Value = 23
$ws.Range("AD10").Value = 6.5
$ws.Range("AI10").Value = 12
$ws.Range("AK10").Value = 26
$ws.Range("AP10").Value = 26
$ws.Range("AW10").Value = 4.5
$ws.Range("BA10").Value = 81
